$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 15:25"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1747781
$ws.Range("C4").Value = 1978
$ws.Range("E4").Value = 1155328
$ws.Range("G4").Value = 90
$ws.Range("H4").Value = 102197

# Row 13 - India
$ws.Range("B13").Value = 159138
$ws.Range("C13").Value = 1052
$ws.Range("D13").Value = 67983
$ws.Range("E13").Value = 86613
$ws.Range("G13").Value = 8
$ws.Range("H13").Value = 4542

# Row 59 - Noruega
$ws.Range("B59").Value = 8406
$ws.Range("C59").Value = 5
$ws.Range("E59").Value = 443

# Row 81 - now Republica de Yibuti (was Grecia; country list resorted and new
# country data inserted, pushing Grecia down to row 82)
$ws.Range("A81").Value = "Republica de Yibuti"
$ws.Range("B81").Value = 2914
$ws.Range("C81").Value = 217
$ws.Range("D81").Value = 1241
$ws.Range("E81").Value = 1653
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 20

# Row 82 - now Grecia (shifted down one row, keeps previous Grecia figures)
$ws.Range("A82").Value = "Grecia"
$ws.Range("B82").Value = 2903
$ws.Range("D82").Value = 1374
$ws.Range("E82").Value = 1356
$ws.Range("H82").Value = 173

# Row 91 - Cuba
$ws.Range("B91").Value = 1983
$ws.Range("C91").Value = 9
$ws.Range("D91").Value = 1734
$ws.Range("E91").Value = 167

# Row 101 - Sri Lanka
$ws.Range("B101").Value = 1503
$ws.Range("C101").Value = 34
$ws.Range("E101").Value = 748
